# Fixed CDS files tab function
# The "SamplesTab" row (row 3) had its Cypher query rewritten so that the
# Tumor column is populated via a COLLECT(DISTINCT ...) aggregation and the
# sample node is matched directly in the MATCH pattern instead of via a
# separate OPTIONAL MATCH (also dropping the unused diagnosis/file matches).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["GECCO OICR: Molecular Pathological Epidemiology of Colorectal Cancer"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(tumor,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

# Row 3 = SamplesTab: update its query text (column B) with the fixed query.
$ws.Range("B3").Value = $newSamplesQuery

# The new query text is shorter, so the row autofits to a smaller height.
$ws.Rows.Item(3).RowHeight = 187.2

# Move the active selection to B4 (FilesTab row), matching the saved view state.
$ws.Range("B4").Select()
